$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (investment_value) so the new
# "company_lei" column becomes column D, shifting investment_value,
# engagement_targets, additional_field_1 and additional_field_2 one
# column to the right (columns E:H).
$ws.Columns.Item(4).Insert()

# Header + data for the new company_lei column.
$ws.Cells.Item(1, 4).Value = "company_lei"
$ws.Cells.Item(2, 4).Value = "JP0000000001"
$ws.Cells.Item(3, 4).Value = "UK0000000002"

# Match the new column's width (no auto "best fit", fixed width).
$ws.Columns.Item(4).ColumnWidth = 13.17

# The worksheet table ("Table7") needs to grow to include the new column
# in the correct position. Rebuild it (Unlist + re-Add) so the column
# order in the table definition matches the new physical column order:
# company_name, company_id, company_isin, company_lei, investment_value,
# engagement_targets, additional_field_1, additional_field_2.
$lo = $ws.ListObjects.Item(1)
$lo.Unlist()

$newLo = $ws.ListObjects.Add(1, $ws.Range("A1:H52"), $null, 1)
$newLo.Name = "Table7"
